# demo.xlsx: translate Sheet2 ("Table 2") headers/labels from Chinese to
# English (mirroring Sheet1's already-English text), resize Sheet2's first
# two columns to fit the new text, and update both sheets' saved selections.

$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Table 2"): Chinese -> English -----------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value2 = "String"
$ws2.Range("B1").Value2 = "Date"
$ws2.Range("C1").Value2 = "Number"

$ws2.Range("A2").Value2  = "Table 2 String 0"
$ws2.Range("A3").Value2  = "Table 2 String 1"
$ws2.Range("A4").Value2  = "Table 2 String 2"
$ws2.Range("A5").Value2  = "Table 2 String 3"
$ws2.Range("A6").Value2  = "Table 2 String 4"
$ws2.Range("A7").Value2  = "Table 2 String 5"
$ws2.Range("A8").Value2  = "Table 2 String 6"
$ws2.Range("A9").Value2  = "Table 2 String 7"
$ws2.Range("A10").Value2 = "Table 2 String 8"
$ws2.Range("A11").Value2 = "Table 2 String 9"

# Widen columns A/B on Sheet2 so the longer English labels fit.
$ws2.Columns.Item(1).ColumnWidth = 17.949776785714285
$ws2.Columns.Item(2).ColumnWidth = 19.949776785714285

# Sheet2 ends up active with A16 selected (just below the data).
$ws2.Activate()
$ws2.Range("A16").Select()

# --- Sheet1: restore it as the active tab with a fresh selection --------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A1:C1").Select()
